# Sort the Cylinder summary rows (A5:G8) ascending by swapping the two rows
# that are out of order (row 5 holds Cylinder=6, row 6 holds Cylinder=4), so
# the final column-A order reads 4, 6, 6, 8. Only the values move - each
# cell keeps its own pre-existing formatting. Afterwards the two rows that
# now share Cylinder=6 (rows 6 and 7) get merged into a single, visually
# grouped cell in column A, the way a sorted/grouped key column looks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap row 5 and row 6 values (columns A-G) - this is the only reordering
# needed to put column A in ascending order (4, 6, 6, 8).
for ($c = 1; $c -le 7; $c++) {
  $v5 = $ws.Cells.Item(5, $c).Value()
  $v6 = $ws.Cells.Item(6, $c).Value()
  $ws.Cells.Item(5, $c).Value = $v6
  $ws.Cells.Item(6, $c).Value = $v5
}

# Merge the now-adjacent "6" cells in column A (rows 6-7) into a single cell.
$ws.Range("A6:A7").Merge()

# Give the merged cell a top-aligned vertical alignment, matching the look of
# a grouped/sorted column header cell.
$ws.Range("A6").VerticalAlignment = -4160

# The lower half of the merge (A7) no longer holds its own value - clear it
# and restore the plain "empty" border styling used elsewhere in the table
# (same formatting as the blank spacer cells in column H).
$ws.Range("A7").ClearContents()
$ws.Range("H7").Copy()
$ws.Range("A7").PasteSpecial(-4122)
